# Updates the cryptos price/volume table to the latest scraped snapshot.
# For cells whose new value is a plain number (e.g. "322.16"), a leading
# apostrophe forces Excel to keep it as text (matching the original
# inlineStr/string cell type) instead of silently converting it to a
# numeric value; Style is then reset to "Normal" so no stray number-format
# style gets attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.268.46"
$ws.Range("E2").Value = "  +2.10%  "
$ws.Range("D3").Value = "2.384.80"
$ws.Range("E3").Value = "  +7.32%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'322.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +11.30%  "
$ws.Range("D6").Value = "'105.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.96%  "
$ws.Range("D7").Value = "'0.657"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.14%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "'0.652"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.16%  "
$ws.Range("D10").Value = "'41.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.31%  "
$ws.Range("D11").Value = "'0.0942"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.46%  "
$ws.Range("D12").Value = "'8.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").Value = "'17.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +16.01%  "
$ws.Range("D14").Value = "'1.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("E15").Value = "  +2.87%  "
$ws.Range("D16").Value = "2.747.93"
$ws.Range("E16").Value = "  +7.43%  "
$ws.Range("D17").Value = "2.397.67"
$ws.Range("E17").Value = "  +8.03%  "
$ws.Range("D18").Value = "43.253.84"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("D19").Value = "'0.0000108"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.21%  "
$ws.Range("D20").Value = "'7.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.95%  "
$ws.Range("D21").Value = "'75.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.65%  "
$ws.Range("D22").Value = "'3.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.12%  "
$ws.Range("D23").Value = "'267.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.21%  "
$ws.Range("D24").Value = "'2.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("D25").Value = "'9.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.62%  "
$ws.Range("D26").Value = "'11.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.94%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "'22.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.84%  "
$ws.Range("D29").Value = "'176.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.14%  "
$ws.Range("D30").Value = "'2.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.57%  "
$ws.Range("B31").Value = "WEMIXToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D31").Value = "'3.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.39%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "'37.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").Value = "'0.0927"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.93%  "
$ws.Range("D34").Value = "'5.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.14%  "
$ws.Range("E35").Value = "  +6.67%  "
$ws.Range("E36").Value = "  -2.09%  "
$ws.Range("D37").Value = "'4.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.36%  "
$ws.Range("D38").Value = "'0.0368"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.93%  "
$ws.Range("D39").Value = "'0.109"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.64%  "
$ws.Range("D40").Value = "'2.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +18.31%  "
$ws.Range("E41").Value = "  +20.61%  "
$ws.Range("D42").Value = "'125.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +25.03%  "
$ws.Range("E43").Value = "  +1.22%  "
$ws.Range("D44").Value = "'69.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.49%  "
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "'12.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.32%  "
$ws.Range("D47").Value = "'9.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +14.23%  "
$ws.Range("D48").Value = "'5.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.40%  "
$ws.Range("D49").Value = "'86.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +55.28%  "
$ws.Range("D50").Value = "'1.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.21%  "
$ws.Range("D51").Value = "1.597.87"
$ws.Range("E51").Value = "  +12.41%  "
